$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.99999998592334682
$ws.Range("A2").Value = 0.99669140212205765
$ws.Range("A3").Value = 0.98551809631677578
$ws.Range("A4").Value = 0.9831055530226056
$ws.Range("A5").Value = 0.97784806428916515
$ws.Range("A6").Value = 0.96654978319956897
$ws.Range("A7").Value = 0.96446939120537689
$ws.Range("A8").Value = 0.96091177114340476
$ws.Range("A9").Value = 0.95933827467407307
$ws.Range("A10").Value = 0.95870137559367397
$ws.Range("A11").Value = 0.95861966911459096
$ws.Range("A12").Value = 0.95872069495384604
$ws.Range("A13").Value = 0.96227675793386558
$ws.Range("A14").Value = 0.96033563166616087
$ws.Range("A15").Value = 0.95907853516180064
$ws.Range("A16").Value = 0.95657254633505717
$ws.Range("A17").Value = 0.95286519078581478
$ws.Range("A18").Value = 0.95175634982148416
$ws.Range("A19").Value = 0.99477653131197519
$ws.Range("A20").Value = 0.97973267653674267
$ws.Range("A21").Value = 0.97638913867856048
$ws.Range("A22").Value = 0.97197379304648079
$ws.Range("A23").Value = 0.97848354975501739
$ws.Range("A24").Value = 0.96546290311157201
$ws.Range("A25").Value = 0.95900597243679464
$ws.Range("A26").Value = 0.96578667046076117
$ws.Range("A27").Value = 0.96497633669635696
$ws.Range("A28").Value = 0.96258699581477658
$ws.Range("A29").Value = 0.96155795352869988
$ws.Range("A30").Value = 0.96200303515866814
$ws.Range("A31").Value = 0.97266550763563142
$ws.Range("A32").Value = 0.97334130959611731
$ws.Range("A33").Value = 0.97282132173732738
